$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures.
# Each entry: cell address -> new value. D-column numeric-looking
# values need NumberFormat forced to Text ("@") first, otherwise Excel
# auto-converts the assigned string into a floating-point number.
$updates = @(
    @{ Cell = 'D2'; Value = '25.821.28'; AsText = $false }
    @{ Cell = 'E2'; Value = '  +0.11%  '; AsText = $false }
    @{ Cell = 'D3'; Value = '1.740.30'; AsText = $false }
    @{ Cell = 'E3'; Value = '  -0.87%  '; AsText = $false }
    @{ Cell = 'E4'; Value = '  +0.03%  '; AsText = $false }
    @{ Cell = 'D5'; Value = '225.40'; AsText = $true }
    @{ Cell = 'E5'; Value = '  -4.77%  '; AsText = $false }
    @{ Cell = 'E6'; Value = '  +0.04%  '; AsText = $false }
    @{ Cell = 'D7'; Value = '0.5177'; AsText = $true }
    @{ Cell = 'E7'; Value = '  +2.58%  '; AsText = $false }
    @{ Cell = 'D8'; Value = '0.2735'; AsText = $true }
    @{ Cell = 'E8'; Value = '  +3.33%  '; AsText = $false }
    @{ Cell = 'D9'; Value = '39.08'; AsText = $true }
    @{ Cell = 'E9'; Value = '  -5.14%  '; AsText = $false }
    @{ Cell = 'D10'; Value = '0.06081'; AsText = $true }
    @{ Cell = 'D11'; Value = '1.738.23'; AsText = $false }
    @{ Cell = 'E11'; Value = '  -0.99%  '; AsText = $false }
    @{ Cell = 'D12'; Value = '0.06996'; AsText = $true }
    @{ Cell = 'E12'; Value = '  +1.12%  '; AsText = $false }
    @{ Cell = 'D13'; Value = '15.14'; AsText = $true }
    @{ Cell = 'E13'; Value = '  -3.21%  '; AsText = $false }
    @{ Cell = 'D14'; Value = '0.6307'; AsText = $true }
    @{ Cell = 'E14'; Value = '  +5.53%  '; AsText = $false }
    @{ Cell = 'D15'; Value = '4.485'; AsText = $true }
    @{ Cell = 'E15'; Value = '  +0.21%  '; AsText = $false }
    @{ Cell = 'D16'; Value = '76.38'; AsText = $true }
    @{ Cell = 'E16'; Value = '  -1.23%  '; AsText = $false }
    @{ Cell = 'E17'; Value = '  +0.05%  '; AsText = $false }
    @{ Cell = 'D18'; Value = '1.001'; AsText = $true }
    @{ Cell = 'E18'; Value = '  +0.06%  '; AsText = $false }
    @{ Cell = 'D19'; Value = '25.838.05'; AsText = $false }
    @{ Cell = 'E19'; Value = '  +0.05%  '; AsText = $false }
    @{ Cell = 'D20'; Value = '11.44'; AsText = $true }
    @{ Cell = 'E20'; Value = '  -1.62%  '; AsText = $false }
    @{ Cell = 'D21'; Value = '0.000006612'; AsText = $true }
    @{ Cell = 'E21'; Value = '  -2.82%  '; AsText = $false }
    @{ Cell = 'D22'; Value = '1.963.07'; AsText = $false }
    @{ Cell = 'E22'; Value = '  -0.65%  '; AsText = $false }
    @{ Cell = 'D23'; Value = '4.070'; AsText = $true }
    @{ Cell = 'E23'; Value = '  -0.02%  '; AsText = $false }
    @{ Cell = 'D24'; Value = '8.424'; AsText = $true }
    @{ Cell = 'E24'; Value = '  +2.27%  '; AsText = $false }
    @{ Cell = 'D25'; Value = '5.079'; AsText = $true }
    @{ Cell = 'E25'; Value = '  -2.21%  '; AsText = $false }
    @{ Cell = 'D26'; Value = '136.74'; AsText = $true }
    @{ Cell = 'E26'; Value = '  -0.58%  '; AsText = $false }
    @{ Cell = 'D27'; Value = '1.499'; AsText = $true }
    @{ Cell = 'E27'; Value = '  +3.28%  '; AsText = $false }
    @{ Cell = 'D28'; Value = '1.815'; AsText = $true }
    @{ Cell = 'E28'; Value = '  -0.53%  '; AsText = $false }
    @{ Cell = 'E29'; Value = '  -0.34%  '; AsText = $false }
    @{ Cell = 'D30'; Value = '102.69'; AsText = $true }
    @{ Cell = 'E30'; Value = '  +0.18%  '; AsText = $false }
    @{ Cell = 'D31'; Value = '0.08298'; AsText = $true }
    @{ Cell = 'E31'; Value = '  +1.43%  '; AsText = $false }
    @{ Cell = 'D32'; Value = '3.611'; AsText = $true }
    @{ Cell = 'E32'; Value = '  -1.35%  '; AsText = $false }
    @{ Cell = 'D33'; Value = '3.374'; AsText = $true }
    @{ Cell = 'E33'; Value = '  +0.08%  '; AsText = $false }
    @{ Cell = 'D34'; Value = '0.04397'; AsText = $true }
    @{ Cell = 'E34'; Value = '  +0.08%  '; AsText = $false }
    @{ Cell = 'D35'; Value = '2.627'; AsText = $true }
    @{ Cell = 'E35'; Value = '  -0.97%  '; AsText = $false }
    @{ Cell = 'D36'; Value = '0.9678'; AsText = $true }
    @{ Cell = 'E36'; Value = '  -2.96%  '; AsText = $false }
    @{ Cell = 'D37'; Value = '0.5959'; AsText = $true }
    @{ Cell = 'E37'; Value = '  -0.83%  '; AsText = $false }
    @{ Cell = 'D38'; Value = '2.677'; AsText = $true }
    @{ Cell = 'E38'; Value = '  -1.68%  '; AsText = $false }
    @{ Cell = 'D39'; Value = '0.01552'; AsText = $true }
    @{ Cell = 'E39'; Value = '  +0.20%  '; AsText = $false }
    @{ Cell = 'D40'; Value = '1.925'; AsText = $true }
    @{ Cell = 'E40'; Value = '  -0.58%  '; AsText = $false }
    @{ Cell = 'D41'; Value = '0.9999'; AsText = $true }
    @{ Cell = 'E41'; Value = '  +0.02%  '; AsText = $false }
    @{ Cell = 'D42'; Value = '101.08'; AsText = $true }
    @{ Cell = 'D43'; Value = '0.3804'; AsText = $true }
    @{ Cell = 'E43'; Value = '  +0.27%  '; AsText = $false }
    @{ Cell = 'D44'; Value = '0.7256'; AsText = $true }
    @{ Cell = 'E44'; Value = '  -1.15%  '; AsText = $false }
    @{ Cell = 'D45'; Value = '4.854'; AsText = $true }
    @{ Cell = 'E45'; Value = '  -1.44%  '; AsText = $false }
    @{ Cell = 'D46'; Value = '0.05483'; AsText = $true }
    @{ Cell = 'E46'; Value = '  -0.16%  '; AsText = $false }
    @{ Cell = 'D47'; Value = '6.151'; AsText = $true }
    @{ Cell = 'E47'; Value = '  +3.87%  '; AsText = $false }
    @{ Cell = 'E48'; Value = '  +0.89%  '; AsText = $false }
    @{ Cell = 'D49'; Value = '29.65'; AsText = $true }
    @{ Cell = 'E49'; Value = '  -0.22%  '; AsText = $false }
    @{ Cell = 'D50'; Value = '51.89'; AsText = $true }
    @{ Cell = 'E50'; Value = '  -0.09%  '; AsText = $false }
    @{ Cell = 'E51'; Value = '  +0.23%  '; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
